$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price (D) and Volume(1h) (E) columns to Text format first so that
# numeric-looking values (e.g. "31.192.86", "253.36", "0.05250") are
# written back as literal text instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '31.192.86'
$ws.Range("E2").Value = '  +2.07%  '

$ws.Range("D3").Value = '1.992.35'
$ws.Range("E3").Value = '  +5.86%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").Value = '0.7772'
$ws.Range("E5").Value = '  +64.27%  '

$ws.Range("D6").Value = '253.36'
$ws.Range("E6").Value = '  +2.93%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").Value = '0.3470'
$ws.Range("E8").Value = '  +20.15%  '

$ws.Range("D9").Value = '27.55'
$ws.Range("E9").Value = '  +23.36%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '44.13'
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.07025'
$ws.Range("E11").Value = '  +7.66%  '

$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.8421'
$ws.Range("E12").Value = '  +10.61%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.08159'
$ws.Range("E13").Value = '  +4.44%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '100.48'
$ws.Range("E14").Value = '  +0.35%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.991.58'
$ws.Range("E15").Value = '  +5.82%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '5.615'
$ws.Range("E16").Value = '  +7.24%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '15.35'
$ws.Range("E17").Value = '  +16.30%  '

$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").Value = '271.69'
$ws.Range("E18").Value = '  -4.56%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '31.194.68'
$ws.Range("E19").Value = '  +2.12%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000008011'
$ws.Range("E20").Value = '  +6.59%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '5.871'
$ws.Range("E21").Value = '  +9.81%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.253.13'
$ws.Range("E22").Value = '  +5.85%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.11%  '

$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '7.047'
$ws.Range("E25").Value = '  +9.51%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.993'
$ws.Range("E26").Value = '  +9.07%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '165.36'
$ws.Range("E27").Value = '  +1.29%  '

$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '0.1441'
$ws.Range("E28").Value = '  +48.77%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '19.81'
$ws.Range("E29").Value = '  +4.01%  '

$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '2.341'
$ws.Range("E30").Value = '  +23.07%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.595'
$ws.Range("E31").Value = '  +6.20%  '

$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = '1.363'
$ws.Range("E32").Value = '  +2.69%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.595'
$ws.Range("E33").Value = '  +8.31%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.425'
$ws.Range("E34").Value = '  +6.08%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.05250'
$ws.Range("E35").Value = '  +8.48%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7818'
$ws.Range("E36").Value = '  +12.05%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '1.218'
$ws.Range("E37").Value = '  +8.04%  '

$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '2.763'
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("D39").Value = '0.02002'
$ws.Range("E39").Value = '  +5.20%  '

$ws.Range("D40").Value = '2.894'
$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("D41").Value = '6.694'
$ws.Range("E41").Value = '  +6.44%  '

$ws.Range("D42").Value = '79.48'
$ws.Range("E42").Value = '  +5.32%  '

$ws.Range("D43").Value = '0.4670'
$ws.Range("E43").Value = '  +10.00%  '

$ws.Range("D44").Value = '2.104'
$ws.Range("E44").Value = '  +6.56%  '

$ws.Range("E45").Value = '  +1.82%  '

$ws.Range("D46").Value = '104.55'
$ws.Range("E46").Value = '  +3.30%  '

$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  +0.20%  '

$ws.Range("D48").Value = '9.971'
$ws.Range("E48").Value = '  +1.14%  '

$ws.Range("D49").Value = '7.644'
$ws.Range("E49").Value = '  +8.95%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '37.31'
$ws.Range("E50").Value = '  +5.93%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.565'
$ws.Range("E51").Value = '  +17.21%  '
